# RPA datasets push 2024-05-01
# Adds a new subscription row (KB / 제일엠앤에스) at row 2 and shifts the
# existing rows down by one, appending what was the last row as new row 10.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("A2").Value = "KB"
$ws.Range("B2").Value = "2024-04-18"
$ws.Range("C2").Value = "제일엠앤에스"
$ws.Range("D2").Value = "KB"
$ws.Range("E2").Value = "KB"
$ws.Range("F2").Value = "2024-04-23"
$ws.Range("G2").Value = "2024-04-30"
$ws.Range("H2").Value = 52800
$ws.Range("I2").Value = 2400000
$ws.Range("J2").Value = 22000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 100
$ws.Range("B2").ClearFormats()
$ws.Range("F2").ClearFormats()
$ws.Range("G2").ClearFormats()

# Row 3
$ws.Range("B3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("A3").Value = "NH"
$ws.Range("B3").Value = "2024-03-14"
$ws.Range("C3").Value = "엔젤로보틱스"
$ws.Range("D3").Value = "NH"
$ws.Range("E3").Value = "NH"
$ws.Range("F3").Value = "2024-03-19"
$ws.Range("G3").Value = "2024-03-26"
$ws.Range("H3").Value = 32000
$ws.Range("I3").Value = 1600000
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100
$ws.Range("B3").ClearFormats()
$ws.Range("F3").ClearFormats()
$ws.Range("G3").ClearFormats()

# Row 4
$ws.Range("B4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("A4").Value = "NH"
$ws.Range("B4").Value = "2024-03-04"
$ws.Range("C4").Value = "오상헬스케어"
$ws.Range("D4").Value = "NH"
$ws.Range("E4").Value = "NH"
$ws.Range("F4").Value = "2024-03-07"
$ws.Range("G4").Value = "2024-03-13"
$ws.Range("H4").Value = 19800
$ws.Range("I4").Value = 990000
$ws.Range("J4").Value = 20000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100
$ws.Range("B4").ClearFormats()
$ws.Range("F4").ClearFormats()
$ws.Range("G4").ClearFormats()

# Row 5
$ws.Range("B5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("A5").Value = "미래"
$ws.Range("B5").Value = "2024-03-25"
$ws.Range("C5").Value = "아이엠비디엑스"
$ws.Range("D5").Value = "미래"
$ws.Range("E5").Value = "미래"
$ws.Range("F5").Value = "2024-03-28"
$ws.Range("G5").Value = "2024-04-03"
$ws.Range("H5").Value = 32500
$ws.Range("I5").Value = 2500000
$ws.Range("J5").Value = 13000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100
$ws.Range("B5").ClearFormats()
$ws.Range("F5").ClearFormats()
$ws.Range("G5").ClearFormats()

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("A6").Value = "신한"
$ws.Range("B6").Value = "2024-04-11"
$ws.Range("C6").Value = "신한제13호스팩"
$ws.Range("D6").Value = "신한"
$ws.Range("E6").Value = "신한"
$ws.Range("F6").Value = "2024-04-15"
$ws.Range("G6").Value = "2024-04-22"
$ws.Range("H6").Value = 6000
$ws.Range("I6").Value = 3000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100
$ws.Range("B6").ClearFormats()
$ws.Range("F6").ClearFormats()
$ws.Range("G6").ClearFormats()

# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("A7").Value = "신한"
$ws.Range("B7").Value = "2024-04-02"
$ws.Range("C7").Value = "신한제12호스팩"
$ws.Range("D7").Value = "신한"
$ws.Range("E7").Value = "신한"
$ws.Range("F7").Value = "2024-04-05"
$ws.Range("G7").Value = "2024-04-15"
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100
$ws.Range("B7").ClearFormats()
$ws.Range("F7").ClearFormats()
$ws.Range("G7").ClearFormats()

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("A8").Value = "하나"
$ws.Range("B8").Value = "2024-04-15"
$ws.Range("C8").Value = "하나33호스팩"
$ws.Range("D8").Value = "하나"
$ws.Range("E8").Value = "하나"
$ws.Range("F8").Value = "2024-04-18"
$ws.Range("G8").Value = "2024-04-24"
$ws.Range("H8").Value = 7000
$ws.Range("I8").Value = 3500000
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100
$ws.Range("B8").ClearFormats()
$ws.Range("F8").ClearFormats()
$ws.Range("G8").ClearFormats()

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("A9").Value = "하나"
$ws.Range("B9").Value = "2024-03-18"
$ws.Range("C9").Value = "하나32호스팩"
$ws.Range("D9").Value = "하나"
$ws.Range("E9").Value = "하나"
$ws.Range("F9").Value = "2024-03-21"
$ws.Range("G9").Value = "2024-03-27"
$ws.Range("H9").Value = 6000
$ws.Range("I9").Value = 3000000
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100
$ws.Range("B9").ClearFormats()
$ws.Range("F9").ClearFormats()
$ws.Range("G9").ClearFormats()

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("A10").Value = "한국"
$ws.Range("B10").Value = "2024-03-12"
$ws.Range("C10").Value = "삼현"
$ws.Range("D10").Value = "한국"
$ws.Range("E10").Value = "한국"
$ws.Range("F10").Value = "2024-03-15"
$ws.Range("G10").Value = "2024-03-21"
$ws.Range("H10").Value = 60000
$ws.Range("I10").Value = 2000000
$ws.Range("J10").Value = 30000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100
$ws.Range("B10").ClearFormats()
$ws.Range("F10").ClearFormats()
$ws.Range("G10").ClearFormats()

$ws.Range("A1:L1").Select()
